$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: both languages have now been handed back ---
$wsOverview.Range("B2").Value = "Handed back"
$wsOverview.Range("C2").Value = "Handed back"

# --- zh-cn sheet: record the handback ---
$wsZhCn.Range("B2").Value = "Handed back"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/286cb87165f7e6956f74dce5bfd5a3772d7d042f/e2e/3a4c0784-b4ee-4359-a849-d72277d74a37.md",
    "",
    "",
    "3a4c0784-b4ee-4359-a849-d72277d74a37.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c2f5cec07b52d67401130139e55d9900b0a054f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/3a4c0784-b4ee-4359-a849-d72277d74a37.020e7190a0749e28a782ba881c3051d83f86f386.zh-cn.xlf",
    "",
    "",
    "3a4c0784-b4ee-4359-a849-d72277d74a37.020e7190a0749e28a782ba881c3051d83f86f386.zh-cn.xlf"
) | Out-Null

$wsZhCn.Range("G2").Value = "2016-01-07 12:50:29"

# --- de-de sheet: record the handback ---
$wsDeDe.Range("B2").Value = "Handed back"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/286cb87165f7e6956f74dce5bfd5a3772d7d042f/e2e/3a4c0784-b4ee-4359-a849-d72277d74a37.md",
    "",
    "",
    "3a4c0784-b4ee-4359-a849-d72277d74a37.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed02c6c0c5df5d257bbe5a7b0eb39d5b9699d083/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/3a4c0784-b4ee-4359-a849-d72277d74a37.020e7190a0749e28a782ba881c3051d83f86f386.de-de.xlf",
    "",
    "",
    "3a4c0784-b4ee-4359-a849-d72277d74a37.020e7190a0749e28a782ba881c3051d83f86f386.de-de.xlf"
) | Out-Null

$wsDeDe.Range("G2").Value = "2016-01-07 12:50:48"

Write-Output "Handback report generated"
